$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 3.639357264516869
$ws.Cells.Item(2, 3).Value = 0.8909090909090909
$ws.Cells.Item(2, 4).Value = 0.8181818181818181
$ws.Cells.Item(3, 2).Value = 1.267300926324749
$ws.Cells.Item(3, 3).Value = 0.8636363636363636
$ws.Cells.Item(3, 4).Value = 0.3
$ws.Cells.Item(4, 2).Value = 1.036731823133914
$ws.Cells.Item(4, 3).Value = 0.2909090909090909
$ws.Cells.Item(4, 4).Value = 0.7181818181818183
$ws.Cells.Item(5, 2).Value = 2.346932117483861
$ws.Cells.Item(5, 3).Value = 0.5
$ws.Cells.Item(5, 4).Value = 0.1363636363636364
$ws.Cells.Item(6, 2).Value = 2.994319709317411
$ws.Cells.Item(6, 3).Value = 0.2363636363636364
$ws.Cells.Item(6, 4).Value = 0.8272727272727273
$ws.Cells.Item(7, 2).Value = 2.628161784490352
$ws.Cells.Item(7, 3).Value = 0.3727272727272727
$ws.Cells.Item(7, 4).Value = 0.2818181818181819
$ws.Cells.Item(8, 2).Value = 1.956867103089487
$ws.Cells.Item(8, 3).Value = 0.5636363636363636
$ws.Cells.Item(8, 4).Value = 0.3090909090909091
$ws.Cells.Item(9, 2).Value = 2.221315499852857
$ws.Cells.Item(9, 3).Value = 0.5727272727272728
$ws.Cells.Item(9, 4).Value = 0.2454545454545455
$ws.Cells.Item(10, 2).Value = 3.210639168481903
$ws.Cells.Item(10, 3).Value = 0.6272727272727272
$ws.Cells.Item(10, 4).Value = 0.6818181818181819
$ws.Cells.Item(11, 2).Value = 0.5490662632377095
$ws.Cells.Item(11, 3).Value = 0.9181818181818181
$ws.Cells.Item(11, 4).Value = 0.9454545454545454
$ws.Cells.Item(12, 2).Value = 1.536370198059503
$ws.Cells.Item(12, 3).Value = 0.2181818181818182
$ws.Cells.Item(12, 4).Value = 0.3727272727272728
$ws.Cells.Item(13, 2).Value = 1.396572294386483
$ws.Cells.Item(13, 3).Value = 0.2454545454545455
$ws.Cells.Item(13, 4).Value = 0.03636363636363638
$ws.Cells.Item(14, 2).Value = 3.316944508690913
$ws.Cells.Item(14, 3).Value = 0.8999999999999999
$ws.Cells.Item(14, 4).Value = 0.03636363636363638
$ws.Cells.Item(15, 2).Value = 3.721807109487508
$ws.Cells.Item(15, 3).Value = 0.5363636363636364
$ws.Cells.Item(15, 4).Value = 0.1909090909090909
$ws.Cells.Item(16, 2).Value = 1.990433615503845
$ws.Cells.Item(16, 3).Value = 0.5363636363636364
$ws.Cells.Item(16, 4).Value = 0.7727272727272727
$ws.Cells.Item(17, 2).Value = 0.7216216206363824
$ws.Cells.Item(17, 3).Value = 0.4181818181818182
$ws.Cells.Item(17, 4).Value = 0.1363636363636364
$ws.Cells.Item(18, 2).Value = 1.649746323567624
$ws.Cells.Item(18, 3).Value = 0.1272727272727273
$ws.Cells.Item(18, 4).Value = 0.3363636363636364
$ws.Cells.Item(19, 2).Value = 2.211936483046026
$ws.Cells.Item(19, 3).Value = 0.3454545454545455
$ws.Cells.Item(19, 4).Value = 0.1090909090909091
$ws.Cells.Item(20, 2).Value = 3.558903765998242
$ws.Cells.Item(20, 3).Value = 0.4636363636363636
$ws.Cells.Item(20, 4).Value = 0.3363636363636364
$ws.Cells.Item(21, 2).Value = 2.052472253384159
$ws.Cells.Item(21, 3).Value = 0.9272727272727272
$ws.Cells.Item(21, 4).Value = 0.5818181818181818
$ws.Cells.Item(22, 2).Value = 1.908738155791507
$ws.Cells.Item(22, 3).Value = 0.8272727272727273
$ws.Cells.Item(22, 4).Value = 0.6818181818181819
$ws.Cells.Item(23, 2).Value = 1.752570065567825
$ws.Cells.Item(23, 3).Value = 0.5545454545454546
$ws.Cells.Item(23, 4).Value = 0.4363636363636364
$ws.Cells.Item(24, 2).Value = 1.356084410781417
$ws.Cells.Item(24, 3).Value = 0.08181818181818182
$ws.Cells.Item(24, 4).Value = 0.04545454545454553
$ws.Cells.Item(25, 2).Value = 1.904446878351196
$ws.Cells.Item(25, 3).Value = 0.6454545454545454
$ws.Cells.Item(25, 4).Value = 0.3090909090909091
$ws.Cells.Item(26, 2).Value = 1.985971858623564
$ws.Cells.Item(26, 3).Value = 0.7545454545454545
$ws.Cells.Item(26, 4).Value = 0.7545454545454545
$ws.Cells.Item(27, 2).Value = 2.616714026348632
$ws.Cells.Item(27, 3).Value = 0.5818181818181818
$ws.Cells.Item(27, 4).Value = 0.1818181818181819
$ws.Cells.Item(28, 2).Value = 2.894649102499058
$ws.Cells.Item(28, 3).Value = 0.9363636363636363
$ws.Cells.Item(28, 4).Value = 0.5272727272727273
$ws.Cells.Item(29, 2).Value = 2.16766252887999
$ws.Cells.Item(29, 3).Value = 0.6727272727272727
$ws.Cells.Item(29, 4).Value = 0.7090909090909091
$ws.Cells.Item(30, 2).Value = 3.153741985409201
$ws.Cells.Item(30, 4).Value = 0.2000000000000001
$ws.Cells.Item(31, 2).Value = 1.928829294193834
$ws.Cells.Item(31, 3).Value = 0.4909090909090909
$ws.Cells.Item(31, 4).Value = 0.2454545454545455
$ws.Cells.Item(32, 2).Value = 1.63125418964535
$ws.Cells.Item(32, 3).Value = 0.3545454545454546
$ws.Cells.Item(32, 4).Value = 0.8636363636363636
$ws.Cells.Item(33, 2).Value = 1.972969181088311
$ws.Cells.Item(33, 3).Value = 0.5909090909090908
$ws.Cells.Item(33, 4).Value = 0.04545454545454553
$ws.Cells.Item(34, 2).Value = 2.2614534407422
$ws.Cells.Item(34, 3).Value = 0.8545454545454545
$ws.Cells.Item(34, 4).Value = 0.3181818181818182
$ws.Cells.Item(35, 2).Value = 3.020645444184856
$ws.Cells.Item(35, 3).Value = 0.09090909090909091
$ws.Cells.Item(35, 4).Value = 0.06363636363636371
$ws.Cells.Item(36, 2).Value = 1.430288620498196
$ws.Cells.Item(36, 3).Value = 0.2545454545454545
$ws.Cells.Item(36, 4).Value = 0.7272727272727273
$ws.Cells.Item(37, 2).Value = 1.73049082626007
$ws.Cells.Item(37, 3).Value = 0.609090909090909
$ws.Cells.Item(37, 4).Value = 0.7181818181818183
$ws.Cells.Item(38, 2).Value = 3.097541125975837
$ws.Cells.Item(38, 3).Value = 0.9545454545454545
$ws.Cells.Item(38, 4).Value = 0.2181818181818183
$ws.Cells.Item(39, 2).Value = 2.664868520100121
$ws.Cells.Item(39, 3).Value = 0.609090909090909
$ws.Cells.Item(39, 4).Value = 0.1090909090909091
$ws.Cells.Item(40, 2).Value = 2.872891472105213
$ws.Cells.Item(40, 3).Value = 0.3090909090909091
$ws.Cells.Item(40, 4).Value = 0.490909090909091
$ws.Cells.Item(41, 2).Value = 2.847233275282147
$ws.Cells.Item(41, 3).Value = 0.4818181818181818
$ws.Cells.Item(41, 4).Value = 0.7090909090909091
$ws.Cells.Item(42, 2).Value = 1.884788789239989
$ws.Cells.Item(42, 3).Value = 0.3636363636363636
$ws.Cells.Item(42, 4).Value = 0.7090909090909091
$ws.Cells.Item(43, 2).Value = 2.052136444473186
$ws.Cells.Item(43, 3).Value = 0.9454545454545454
$ws.Cells.Item(43, 4).Value = 0.1000000000000001
$ws.Cells.Item(44, 2).Value = 1.670744220891262
$ws.Cells.Item(44, 3).Value = 0.05454545454545454
$ws.Cells.Item(44, 4).Value = 0.7181818181818183
$ws.Cells.Item(45, 2).Value = 0.7968207385319868
$ws.Cells.Item(45, 3).Value = 0.5181818181818182
$ws.Cells.Item(45, 4).Value = 0.7181818181818183
$ws.Cells.Item(46, 2).Value = 2.414704811449499
$ws.Cells.Item(46, 3).Value = 0.4363636363636363
$ws.Cells.Item(46, 4).Value = 0.5454545454545454
$ws.Cells.Item(47, 2).Value = 3.637041109218055
$ws.Cells.Item(47, 3).Value = 0.5727272727272728
$ws.Cells.Item(47, 4).Value = 0.2181818181818183
$ws.Cells.Item(48, 2).Value = 1.799200636417428
$ws.Cells.Item(48, 3).Value = 0.7454545454545454
$ws.Cells.Item(48, 4).Value = 0.04545454545454553
$ws.Cells.Item(49, 2).Value = 2.617465268312308
$ws.Cells.Item(49, 3).Value = 0.1363636363636364
$ws.Cells.Item(49, 4).Value = 0.2727272727272727
$ws.Cells.Item(50, 2).Value = 2.549483547470375
$ws.Cells.Item(50, 3).Value = 0.4545454545454545
$ws.Cells.Item(50, 4).Value = 0.8
$ws.Cells.Item(51, 2).Value = 2.166045691038782
$ws.Cells.Item(51, 3).Value = 0.8545454545454545
$ws.Cells.Item(51, 4).Value = 0.5636363636363637